# Fill in row 13 of Sheet1 (the 17/10/2025 -> 24/10/2025 match week's stats)
# which was previously a set of blank placeholder cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "24/10/2025"
$ws.Range("B13").Value = "Motor Lublin"
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = "Widzew Lodz"
$ws.Range("F13").Value = "W"
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 2.65
$ws.Range("L13").Value = 0.51
$ws.Range("M13").Value = 16
$ws.Range("N13").Value = 13
$ws.Range("O13").Value = 5
$ws.Range("P13").Value = 4
